$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 175, shifting the existing rows 175:219 down to 176:220
$ws.Rows("175:175").Insert()

# Populate the newly inserted row 175 with the new weekly record
$ws.Cells.Item(175, 1).Value = 5
$ws.Cells.Item(175, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(175, 3).Value = "Maule"
$ws.Cells.Item(175, 4).Value = 44508
$ws.Cells.Item(175, 5).Value = 7
$ws.Cells.Item(175, 6).Value = 100112023
$ws.Cells.Item(175, 7).Value = "Brócoli"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 5000
$ws.Cells.Item(175, 11).Value = 500
$ws.Cells.Item(175, 12).Value = 500
$ws.Cells.Item(175, 13).Value = 500
$ws.Cells.Item(175, 14).Value = "`$/unidad"
$ws.Cells.Item(175, 15).Value = "Región del Maule"
$ws.Cells.Item(175, 16).Value = 500
$ws.Cells.Item(175, 17).Value = 1
$ws.Cells.Item(175, 18).Value = "Hortaliza"
